$d = $word.ActiveDocument

# Locate the "O - Open/Closed" paragraph; the new "L" paragraph is inserted
# right after it (and before the pre-existing trailing blank paragraph).
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*Abierto/Cerrado*") {
        $anchorPara = $candidate
        break
    }
}

if ($anchorPara -eq $null) {
    throw "Could not locate the Open/Closed (O) paragraph to anchor the insertion."
}

$anchorIndex = $anchorPara.Index

# Create a new paragraph right after the O paragraph; it inherits the O
# paragraph's mark formatting (justify both, Arial 12pt).
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($anchorIndex + 1)
$newRange = $newPara.Range

$boldText = "L - Principio de Sustitución de Liskov:"
$restText = " Los objetos deben poder ser sustituidos por instancias de sus subtipos sin alterar el comportamiento del programa."

$startPos = $newRange.Start
$newRange.Text = $boldText + $restText

$titleRange = $d.Range($startPos, $startPos + $boldText.Length)
$titleRange.Font.Bold = 1

# Add a second, empty paragraph after the new "L" paragraph (matching the
# blank paragraph that follows it in the target document), before the
# document's pre-existing trailing empty paragraph.
$newPara = $d.Paragraphs($anchorIndex + 1)
$newPara.Range.InsertParagraphAfter()

Write-Output "Inserted Liskov paragraph after index $anchorIndex; document now has $($d.Paragraphs.Count) paragraphs."
